$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple +1 increments on column B for rows 21-27 and row 32
$ws.Range("B21").Value = 87164
$ws.Range("B22").Value = 87048
$ws.Range("B23").Value = 90539
$ws.Range("B24").Value = 87199
$ws.Range("B25").Value = 87072
$ws.Range("B26").Value = 87224
$ws.Range("B27").Value = 87211
$ws.Range("B32").Value = 87095

# Rows 28 and 29 swap their taxon data (A,B,E,F,G,H,Q,R), with B also +1
$ws.Range("A28").Value = 130807444
$ws.Range("B28").Value = 91222
$ws.Range("E28").Value = 4188
$ws.Range("F28").Value = "Fransig jordstjärna"
$ws.Range("G28").Value = "Geastrum fimbriatum"
$ws.Range("H28").Value = "Fr.:Pers."
$ws.Range("Q28").Value = 704536
$ws.Range("R28").Value = 6361615

$ws.Range("A29").Value = 130807436
$ws.Range("B29").Value = 87211
$ws.Range("E29").Value = 3674
$ws.Range("F29").Value = "Anisspindling"
$ws.Range("G29").Value = "Cortinarius odorifer"
$ws.Range("H29").Value = "Britzelm."
$ws.Range("Q29").Value = 704389
$ws.Range("R29").Value = 6361480

# Rows 30 and 31 swap their taxon data (A,B,E,F,G,H,Q,R), with B also +1
$ws.Range("A30").Value = 130807434
$ws.Range("B30").Value = 87095
$ws.Range("E30").Value = 424
$ws.Range("F30").Value = "Svartgrön spindling"
$ws.Range("G30").Value = "Cortinarius atrovirens"
$ws.Range("H30").Value = "Kalchbr."
$ws.Range("Q30").Value = 704395
$ws.Range("R30").Value = 6361502

$ws.Range("A31").Value = 130807439
$ws.Range("B31").Value = 87199
$ws.Range("E31").Value = 6003296
$ws.Range("F31").Value = "Stor odörspindling"
$ws.Range("G31").Value = "Cortinarius mussivus"
$ws.Range("H31").Value = "(Fr.) Melot"
$ws.Range("Q31").Value = 704409
$ws.Range("R31").Value = 6361473
